# "Minor Ruby script clean up and re-generate reports"
#
# The report-generating script found one additional "rb" (Ruby) source
# file since the last export, so the "rb" row in the File type / Count
# table goes from 8 -> 9. The two charts on the sheet ("Whole repository
# count of files grouped by type") are built directly from this table
# (range 'Chart Report'!$A$2:$B$25), so they pick up the refreshed count
# the next time the report/workbook is regenerated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# File type "rb" lives on row 8 (A8 = "rb", B8 = count). Bump the count
# from 8 to 9.
$ws.Range("B8").Value = 9
